$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; this shifts all existing champion
# rows down by one (row 1 -> row 2, ... row 60 -> row 61) and keeps their
# existing values (name / attributes / cost) unchanged.
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "champs"
$ws.Range("B1").Value = "attributes"
$ws.Range("C1").Value = "cost"

# Restore the view: top-left at A1, active cell / selection at C1.
$ws.Range("C1").Select()

Write-Output "header row added"
